function Set-RowValues($ws, $row, $vals) {
    $n = $vals.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    $startCell = $ws.Cells.Item($row, 1)
    $endCell = $ws.Cells.Item($row, $n)
    $range = $ws.Range($startCell, $endCell)
    $range.Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-RowValues $ws 2 @("ECs", "Apoe", "Vldlr", "ECs", 3, 1, 156.656447, 469.969341, 0.0671576211124673, 0.0671576211124673, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 17.42406109875933, 156.816549888834, 0.0005377552588021451, 0.000537755258802145)
Set-RowValues $ws 3 @("ECs", "Apoe", "Vldlr", "FAPs", 3, 1, 156.656447, 469.969341, 0.0671576211124673, 0.0671576211124673, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 1852.498203444916, 16672.48383100424, 0.05717327580393773, 0.05717327580393772)
Set-RowValues $ws 4 @("ECs", "Apoe", "Vldlr", "M2", 3, 1, 156.656447, 469.969341, 0.0671576211124673, 0.0671576211124673, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 5.943127848654665, 53.488150637892, 0.0001834215477225967, 0.0001834215477225967)
Set-RowValues $ws 5 @("ECs", "Apoe", "Vldlr", "sCs", 3, 1, 156.656447, 469.969341, 0.0671576211124673, 0.0671576211124673, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 300.140279997558, 2701.262519978022, 0.009263168502004827, 0.009263168502004827)
Set-RowValues $ws 6 @("FAPs", "Apoe", "Vldlr", "ECs", 3, 1, 42.300692, 126.902076, 0.01813403725498241, 0.01813403725498241, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 4.704880367469333, 42.343923307224, 0.0001452057672032472, 0.0001452057672032471)
Set-RowValues $ws 7 @("FAPs", "Apoe", "Vldlr", "FAPs", 3, 1, 42.300692, 126.902076, 0.01813403725498241, 0.01813403725498241, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 500.215327457776, 4501.937947119984, 0.01543804405581484, 0.01543804405581484)
Set-RowValues $ws 8 @("FAPs", "Apoe", "Vldlr", "M2", 3, 1, 42.300692, 126.902076, 0.01813403725498241, 0.01813403725498241, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 1.604775452634666, 14.442979073712, 0.00004952785885905394, 0.00004952785885905394)
Set-RowValues $ws 9 @("FAPs", "Apoe", "Vldlr", "sCs", 3, 1, 42.300692, 126.902076, 0.01813403725498241, 0.01813403725498241, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 81.044488012488, 729.4003921123921, 0.002501259573105265, 0.002501259573105265)
Set-RowValues $ws 10 @("M2", "Apoe", "Vldlr", "ECs", 3, 1, 2110.189616, 6330.568848, 0.9046248489651427, 0.9046248489651426, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 234.7051366430614, 2112.346229787552, 0.007243656962765659, 0.007243656962765656)
Set-RowValues $ws 11 @("M2", "Apoe", "Vldlr", "FAPs", 3, 1, 2110.189616, 6330.568848, 0.9046248489651427, 0.9046248489651426, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 24953.47333243245, 224581.259991892, 0.7701339793195582, 0.770133979319558)
Set-RowValues $ws 12 @("M2", "Apoe", "Vldlr", "M2", 3, 1, 2110.189616, 6330.568848, 0.9046248489651427, 0.9046248489651426, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 80.05496685873067, 720.4947017285759, 0.002470720182712123, 0.002470720182712123)
Set-RowValues $ws 13 @("M2", "Apoe", "Vldlr", "sCs", 3, 1, 2110.189616, 6330.568848, 0.9046248489651427, 0.9046248489651426, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 4042.941827949025, 36386.47645154122, 0.1247764925001067, 0.1247764925001067)
Set-RowValues $ws 14 @("sCs", "Apoe", "Vldlr", "ECs", 3, 1, 23.52144266666667, 70.564328, 0.01008349266740757, 0.01008349266740757, 1, 0.3333333333333333, 0.1112246666666667, 0.333674, 0.008007360146089436, 0.008007360146089434, 2.616164620119112, 23.545481581072, 0.00008074215731838444, 0.00008074215731838441)
Set-RowValues $ws 15 @("sCs", "Apoe", "Vldlr", "FAPs", 3, 1, 23.52144266666667, 70.564328, 0.01008349266740757, 0.01008349266740757, 3, 1, 11.825228, 35.475684, 0.851329675722, 0.8513296757219999, 278.1464224222614, 2503.317801800352, 0.008584376542689252, 0.00858437654268925)
Set-RowValues $ws 16 @("sCs", "Apoe", "Vldlr", "M2", 3, 1, 23.52144266666667, 70.564328, 0.01008349266740757, 0.01008349266740757, 1, 0.3333333333333333, 0.03793733333333333, 0.113812, 0.002731209722503793, 0.002731209722503793, 0.8923408109262222, 8.031067298336, 0.00002754013321001927, 0.00002754013321001926)
Set-RowValues $ws 17 @("sCs", "Apoe", "Vldlr", "sCs", 3, 1, 23.52144266666667, 70.564328, 0.01008349266740757, 0.01008349266740757, 3, 1, 1.915914, 5.747742000000001, 0.1379317544094067, 0.1379317544094067, 45.06506130526401, 405.5855517473761, 0.001390833834189915, 0.001390833834189914)

